$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# ---------------------------------------------------------------------------
# 1) "Gesti" + <bookmark _GoBack/> + "onnaire :" -> single run "Gestionnaire :"
#    (the bookmark that used to sit inside this word is removed from here;
#     it reappears below, inside the keywords placeholder.)
# ---------------------------------------------------------------------------
$gestTarget = "Gestionnaire" + $nbsp + ":"
$d.Content.Find.Execute($gestTarget, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $gestTarget, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Keywords placeholder: "{-w:p keywords}{.} : {/keywords}"
#    becomes "{-w:p keywords}{key<bookmark _GoBack/>word} : {definition}{/keywords}"
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$kwStart = $full.IndexOf("{-w:p keywords}{.}")
$dotIdx = $full.IndexOf(".", $kwStart)

# Temporary bookmarks pin down the run boundaries on either side of the "."
# placeholder so that replacing its text does not let the surrounding,
# identically-formatted runs coalesce into it.
$d.Bookmarks.Add("zz_before_dot", $d.Range($dotIdx, $dotIdx)) | Out-Null
$d.Bookmarks.Add("zz_after_dot", $d.Range($dotIdx + 1, $dotIdx + 1)) | Out-Null

$d.Range($dotIdx, $dotIdx + 1).Text = "keyword"

# Split "keyword" into "key" / "word" runs, with the real _GoBack bookmark
# placed right after "word" (mirrors the target markup ordering).
$d.Bookmarks.Add("zz_key_word", $d.Range($dotIdx + 3, $dotIdx + 3)) | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($dotIdx + 7, $dotIdx + 7)) | Out-Null
$d.Bookmarks("zz_key_word").Delete()

$d.Bookmarks("zz_before_dot").Delete()
$d.Bookmarks("zz_after_dot").Delete()

# ---------------------------------------------------------------------------
# 3) Add the "{definition}" run right after the existing single space and
#    before "{/keywords}", keeping that space run untouched/separate.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$kwStart = $full.IndexOf("{-w:p keywords}")
$closeIdx = $full.IndexOf("{/keywords}", $kwStart)
$spaceIdx = $closeIdx - 1

$d.Bookmarks.Add("zz_before_space", $d.Range($spaceIdx, $spaceIdx)) | Out-Null
$d.Bookmarks.Add("zz_after_space", $d.Range($spaceIdx + 1, $spaceIdx + 1)) | Out-Null

$d.Range($spaceIdx, $spaceIdx + 1).Text = " {definition}"

$d.Bookmarks.Add("zz_before_def", $d.Range($spaceIdx + 1, $spaceIdx + 1)) | Out-Null

$d.Bookmarks("zz_before_space").Delete()
$d.Bookmarks("zz_after_space").Delete()

# Re-stamp the {definition} run's text: the in-place split above leaves a
# stray xml:space="preserve" inherited from the space it was split off from;
# cycling the text clears it since the run no longer starts/ends with a space.
$defRange = $d.Range($spaceIdx + 1, $spaceIdx + 1 + 12)
$defRange.Text = "............"
$defRange2 = $d.Range($spaceIdx + 1, $spaceIdx + 1 + 12)
$defRange2.Text = "{definition}"

$d.Bookmarks("zz_before_def").Delete()
